$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("n13", "n13_IMG_3175.jpeg", "True", "no_meltpatch", "negative"),
    @("n14", "n14_IMG_3180.jpeg", "True", "no_meltpatch", "negative"),
    @("n15", "n15_IMG_3174.jpeg", "True", "no_meltpatch", "negative"),
    @("n16", "n16_IMG_3177.jpeg", "True", "no_meltpatch", "negative")
)

$startRow = 14
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = "'" + $rowData[0]
    $ws.Cells.Item($row, 2).Value = "'" + $rowData[1]
    $ws.Cells.Item($row, 3).Value = "'" + $rowData[2]
    $ws.Cells.Item($row, 4).Value = "'" + $rowData[3]
    $ws.Cells.Item($row, 5).Value = "'" + $rowData[4]
}
